# Insert 3 new weekly-report rows for "Tuna" right above the old row 235
# block (Región Metropolitana, week of 2022-02-14 / serial 44617), pushing
# the rest of the data (old rows 235-280) down to 238-283.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at row 235 (each Insert() pushes current row 235
# and everything below it down by one).
$ws.Rows.Item(235).Insert()
$ws.Rows.Item(235).Insert()
$ws.Rows.Item(235).Insert()

# New data for the three inserted rows.
$newRows = @(
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44617, 13, "Fruta", 100107, "Otros", 100107011, "Tuna", "Sin especificar", "Especial", 100, 17000, 17000, 17000, "`$/caja 18 kilos", "Región Metropolitana", 944, 18),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44617, 13, "Fruta", 100107, "Otros", 100107011, "Tuna", "Sin especificar", "Primera", 175, 15000, 15000, 15000, "`$/caja 18 kilos", "Región Metropolitana", 833, 18),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44617, 13, "Fruta", 100107, "Otros", 100107011, "Tuna", "Sin especificar", "Segunda", 175, 11000, 11000, 11000, "`$/caja 18 kilos", "Región Metropolitana", 611, 18)
)

$startRow = 235
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $rowValues[$c - 1]
    }
    # Column D (Fecha) uses the date number format already present on the sheet.
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(239, 4).NumberFormat
}
